$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New 20-minute trade row (row 4) - close price (BuyPrice/SellPrice etc.) can be
# left blank when trader can't locate data from Yahoo.
$ws.Range("A4").Value = 10044.120000000001
$ws.Range("B4").Value = 10119
$ws.Range("C4").Value = 20.3
$ws.Range("D4").Value = 20.149999999999999
$ws.Range("E4").Value = $false
$ws.Range("F4").Value = -0.74
$ws.Range("G4").Value = 42608.640451388892
$ws.Range("H4").Value = $false

# Match the date/time number formatting already used in column G (same as G1/G3).
$ws.Range("G4").NumberFormat = "m/d/yy h:mm"

# Column A needs to grow very slightly to fit the new, slightly wider value.
$ws.Columns.Item(1).ColumnWidth = 8.14
